# This script rotates the species-record data found in rows 3-6 of the
# active sheet. The data that lived in row 3 moves to row 4, row 4's data
# moves to row 5, row 5's data moves to row 6, and row 6's data wraps
# around back to row 3 (i.e. a cyclic "rotate down by one, row 6 -> row 3").
#
# Only the columns that actually differ between these rows are touched
# (A, B, D, E, F, G, H, I, J, P, Q, R); every other column already holds
# identical values across rows 3-6, so leaving them untouched is correct.
#
# We use Range.Copy (cell-to-cell) rather than reading/writing .Value as a
# string, because some of the affected cells (column I) hold a text value
# that merely looks numeric ("3"); round-tripping it through a string
# assignment would make Excel auto-convert it into a real number and
# change its underlying type. Range.Copy preserves the original cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","I","J","P","Q","R")
$sourceRows = @(3,4,5,6)

# New row N receives the values that used to live in the previous row
# (3<-6, 4<-3, 5<-4, 6<-5).
$mapping = @{ 3 = 6; 4 = 3; 5 = 4; 6 = 5 }

# A scratch area well outside the used range, used as a temporary staging
# spot so that overlapping row-to-row copies don't clobber data that is
# still needed for a later step. Cleared again at the end of the script.
$stagingRow = 500

# Stage 1: snapshot every needed source cell into the scratch row, cell by
# cell, preserving each cell's original value/type via Copy.
foreach ($r in $sourceRows) {
    foreach ($col in $cols) {
        $srcCell = $ws.Range($col + $r)
        $stageCell = $ws.Range($col + $stagingRow + $r)
        if ([string]$srcCell.Value2 -eq "") {
            $stageCell.ClearContents()
        } else {
            $srcCell.Copy($stageCell)
        }
    }
}

# Stage 2: write the staged values into their new destination rows.
foreach ($destRow in $sourceRows) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $stageCell = $ws.Range($col + $stagingRow + $srcRow)
        $destCell = $ws.Range($col + $destRow)
        if ([string]$stageCell.Value2 -eq "") {
            $destCell.ClearContents()
        } else {
            $stageCell.Copy($destCell)
        }
    }
}

# Stage 3: clean up the scratch area so it leaves no trace in the saved file.
foreach ($r in $sourceRows) {
    foreach ($col in $cols) {
        $ws.Range($col + $stagingRow + $r).ClearContents()
    }
}
